# Add a new sentence ("... Diagrams are listed with .png. ") to the end
# of the single paragraph, right after the existing " database. " run
# and before the _GoBack bookmark. The new sentence needs "png" wrapped
# in <w:proofErr w:type="spellStart"/>...<w:proofErr w:type="spellEnd"/>
# (matching the existing spell-check markup already used for
# "tt_persons" in this document), which means it has to land in its own
# run, flanked by two more runs. A plain Find/Replace collapses
# everything into a single run and there is no COM call that inserts a
# bare <w:proofErr/> element next to an existing run, so the paragraph's
# whole OOXML is rebuilt and pushed back in one shot via Range.InsertXML
# (InsertXML replaces the exact range it is called on; calling it on a
# range that spans an entire paragraph keeps it a single paragraph
# instead of splitting into a new one).

$d = $word.ActiveDocument

# Locate the paragraph that ends with " database. " (the anchor point
# for the new sentence) without hard-coding paragraph indices.
$anchor = $d.Content
$anchor.Find.Execute(" database. ", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0)
if (-not $anchor.Find.Found) {
    throw "anchor text ' database. ' not found"
}

$para = $anchor.Paragraphs(1)
$prange = $para.Range

# Pull out the bookmark (_GoBack) XML fragment that must stay at the
# very end of the paragraph, and the text that precedes the insertion
# point, straight from the live object model (no hard-coded copy of the
# surrounding sentences).
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$preText = $d.Range($prange.Start, $anchor.End).Text

function Esc([string]$s) {
    $s.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
}

$newSentenceXml = (
    "<w:r><w:t xml:space=`"preserve`">$(Esc ' Diagrams are listed with .')</w:t></w:r>" +
    '<w:proofErr w:type="spellStart"/>' +
    "<w:r><w:t>$(Esc 'png')</w:t></w:r>" +
    '<w:proofErr w:type="spellEnd"/>' +
    "<w:r><w:t xml:space=`"preserve`">$(Esc '. ')</w:t></w:r>"
)

$paraXml = "<w:p $wNs>" +
           "<w:r><w:t xml:space=`"preserve`">$(Esc $preText)</w:t></w:r>" +
           $newSentenceXml +
           '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
           '<w:bookmarkEnd w:id="0"/>' +
           "</w:p>"

$prange.InsertXML($paraXml)
